# Apply the "Chapter 3" workbook edits:
#  - G2: literal text "=D5+F5*F4" (quote-prefixed, so it is stored as text,
#        not evaluated as a formula) - this documents the old G5 formula.
#  - G5: replace the old formula (=D5+F5*F4) with a directly computed
#        value formula (=41+50.73/60), which is the "resolved" sextant
#        reading for this problem.
#  - Move the view/selection to reflect where the user ended up working
#    (cell E13), and let Excel re-settle the scrolled-to cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2 gets the literal (quote-prefixed) text of the prior G5 formula, as a
# little note of what G5 used to contain. The leading apostrophe forces
# Excel to store this as text (quotePrefix style) instead of parsing it
# as a formula.
$ws.Range("G2").Value = "'=D5+F5*F4"

# G5 now computes the sextant altitude directly instead of deriving it
# from D5/F5/F4.
$ws.Range("G5").Formula = "=41+50.73/60"

# Recalculate so all of the dependent cells (E13:L14, etc.) pick up
# the new G5 value.
$excel.CalculateFullRebuild()

# Leave the selection on E13, matching where the user was working.
$ws.Range("E13").Select()
